# Weekly price update: insert two new rows of data (Ajo - Chino - Primera,
# date 2021-11-22 / serial 44522) right after the existing row 370, pushing
# all the subsequent historical rows down by two positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("371:372").Insert()

# New row 371: $/caja 10 kilos
$ws.Range("A371").Value = 10
$ws.Range("B371").Value = "Vega Modelo de Temuco"
$ws.Range("C371").Value = "La Araucanía"
$ws.Range("D371").Value = 44522
$ws.Range("E371").Value = 9
$ws.Range("F371").Value = 100112003
$ws.Range("G371").Value = "Ajo"
$ws.Range("H371").Value = "Chino"
$ws.Range("I371").Value = "Primera"
$ws.Range("J371").Value = 400
$ws.Range("K371").Value = 20000
$ws.Range("L371").Value = 20000
$ws.Range("M371").Value = 20000
$ws.Range("N371").Value = "$/caja 10 kilos"
$ws.Range("O371").Value = "China"
$ws.Range("P371").Value = 2000
$ws.Range("Q371").Value = 10
$ws.Range("R371").Value = "Hortaliza"

# New row 372: $/malla 10 kilos
$ws.Range("A372").Value = 10
$ws.Range("B372").Value = "Vega Modelo de Temuco"
$ws.Range("C372").Value = "La Araucanía"
$ws.Range("D372").Value = 44522
$ws.Range("E372").Value = 9
$ws.Range("F372").Value = 100112003
$ws.Range("G372").Value = "Ajo"
$ws.Range("H372").Value = "Chino"
$ws.Range("I372").Value = "Primera"
$ws.Range("J372").Value = 300
$ws.Range("K372").Value = 20000
$ws.Range("L372").Value = 20000
$ws.Range("M372").Value = 20000
$ws.Range("N372").Value = "$/malla 10 kilos"
$ws.Range("O372").Value = "China"
$ws.Range("P372").Value = 2000
$ws.Range("Q372").Value = 10
$ws.Range("R372").Value = "Hortaliza"
